$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update prices in D33:D36
$ws.Range("D33").Value = 1370
$ws.Range("D34").Value = 1771
$ws.Range("D35").Value = 1980
$ws.Range("D36").Value = 2210
